$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD2:AD53").Value = 86
$ws.Range("AE2:AE53").Value = 76
$ws.Range("AF2:AF53").Value = 0

Write-Output "done"
